$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "product" -> "fsdgsdgdsgh"
$ws.Range("B1").Value = "fsdgsdgdsgh"

# Row 2: B2/C2/D2 used to hold "Motor Insurance" / "Scooters" /
# "Bike 150 CC To 600 CC" with the highlighted "category" style; they now
# just repeat the company template value "hdfc" using the plain style that
# the rest of the row already uses.
$ws.Range("F2").Copy()
$ws.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Value = "hdfc"
$ws.Range("C2").Value = "hdfc"
$ws.Range("D2").Value = "hdfc"

# Move the active selection to D2
$ws.Range("D2").Select()
